# Regression tool updated: new regression tests added to the NewSemTests sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewSemTests")

# New cell on row 4 (col I) documenting the runtime-bug detection column.
$ws.Range("I4").Value = "this test detected a runtime bug"

# Make room for six new rows starting at row 33; this pushes the existing
# "lvalues_runtimeError" row (33) down to row 39, and the four
# "TwoMachines" rows (36-39) down to rows 42-45.
for ($i = 0; $i -lt 6; $i++) {
    $ws.Rows.Item(33).Insert()
}

# --- New rows 32-34: additional SEM_OneMachine regression tests ---------
$ws.Range("A32").Value = "SEM_OneMachine_30\DeferIgnore4"
$ws.Range("B32").Value = "P semantics test, one machine: ""defer"" semantics and the state stack"
$ws.Range("C32").Value = "No"
$ws.Range("D32").Value = "Yes"
$ws.Range("F32").Value = "Yes"
$ws.Range("G32").Value = """unhandled event"" exception wrt deferred event"

$ws.Range("A33").Value = "SEM_OneMachine_31\RaisedHalt"
$ws.Range("B33").Value = "P semantics test: one machine, ""halt"" is raised and unhandled"
$ws.Range("C33").Value = "No"
$ws.Range("D33").Value = "No"
$ws.Range("E33").Value = "Yes"
$ws.Range("F33").Value = "Yes"

$ws.Range("A34").Value = "SEM_OneMachine_32\RaisedHaltHandled"
$ws.Range("B34").Value = "P semantics test: one machine, ""halt"" is raised and handled"
$ws.Range("C34").Value = "No"
$ws.Range("D34").Value = "Yes"
$ws.Range("F34").Value = "Yes"

# --- Row 45 (was row 39 before the insert): "Yes" correction -------------
$ws.Range("D45").Value = "Yes"

# --- New rows 46-48: additional SEM_TwoMachines regression tests ---------
$ws.Range("A46").Value = "SEM_TwoMachines_5\RaisedHalt"
$ws.Range("B46").Value = "P semantics test: two machines, machine is halted with ""raise halt"" (unhandled)"
$ws.Range("C46").Value = "No"
$ws.Range("D46").Value = "No"
$ws.Range("E46").Value = "Yes"
$ws.Range("F46").Value = "Yes"

$ws.Range("A47").Value = "SEM_TwoMachines_6\RaisedHaltHandled"
$ws.Range("B47").Value = "P semantics test: two machines, machine is halted with ""raise halt"" (handled)"
$ws.Range("C47").Value = "No"
$ws.Range("D47").Value = "Yes"
$ws.Range("F47").Value = "Yes"

$ws.Range("A48").Value = "SEM_TwoMachines_7\RaisedHalt_bugFound"
$ws.Range("B48").Value = "P semantics test: two machines, machine is halted with ""raise halt"" (unhandled)"
$ws.Range("C48").Value = "No"
$ws.Range("D48").Value = "No"
$ws.Range("E48").Value = "Yes"
$ws.Range("F48").Value = "Yes"
$ws.Range("I48").Value = "this test found a bug"

# --- View state: scroll down and select row 37 ---------------------------
$ws.Activate()
$ws.Range("A37:XFD37").Select()
$excel.ActiveWindow.ScrollRow = 12

Write-Output "edit complete"
